$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1885723333333333
$ws.Range("H2").Value = 0.565717
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01499333333333333
$ws.Range("N2").Value = 0.04498
$ws.Range("Q2").Value = 0.002827327851111111
$ws.Range("R2").Value = 0.02544595066
